$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value in E3 from -8000 to -80000
$ws.Range("E3").Value = -80000

# Update the active selection to E3
$ws.Range("E3").Select()
